$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the new columns needed for the reshuffled KYC import layout.
#    Original layout: A Investor | B Full Name | C PAN | D Address |
#                      E KYC Type | F Residency | G Bank Account |
#                      H IFSC Code | I Verified | J Update Only |
#                      K Send Kyc Form To User
# ---------------------------------------------------------------------------

# Insert "Date of Birth" before PAN (old C)
$ws.Columns("C").Insert()

# Insert "Correspondence Address" right after Address (now column E)
$ws.Columns("F").Insert()

# "Bank Account" (now column I) becomes "Bank Name"; insert three columns
# after it for Branch Name, Bank Account Number and Account Type.
$ws.Columns("J:L").Insert()

# ---------------------------------------------------------------------------
# 2. Headers (row 1)
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "Date of Birth"
$ws.Range("F1").Value = "Correspondence Address"
$ws.Range("I1").Value = "Bank Name"
$ws.Range("J1").Value = "Branch Name"
$ws.Range("K1").Value = "Bank Account Number"
$ws.Range("L1").Value = "Account Type"

# ---------------------------------------------------------------------------
# 3. Data rows
# ---------------------------------------------------------------------------

# Row 2 (Investor 1) - the old "Bank Account" number (from the pre-insert
# column G) now sits in I2 after the column inserts; move it over to the
# new "Bank Account Number" column (K) and replace I2 with the Bank Name.
$ws.Range("K2").Value = $ws.Range("I2").Value()
$ws.Range("C2").Value = 27478
$ws.Range("I2").Value = "BOB"
$ws.Range("J2").Value = "Xyz"
$ws.Range("L2").Value = "Savings"

# Row 3 (Investor 2)
$ws.Range("K3").Value = $ws.Range("I3").Value()
$ws.Range("C3").Value = 15772
$ws.Range("I3").Value = "HDFC"
$ws.Range("J3").Value = "Abc"
$ws.Range("L3").Value = "Current"

# Apply the date number format to the Date of Birth cells. Format C2
# directly, then copy/paste its formatting onto C3 so both cells share
# the very same style entry (matches the single new cellXfs record).
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Column widths (cosmetic re-sizing to match the new columns; the values
#    below are chosen so the stored OOXML column width lands on the target)
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 9.857142857142858
$ws.Columns("I:J").ColumnWidth = 10.285714285714286
$ws.Columns("K:L").ColumnWidth = 18.0

# ---------------------------------------------------------------------------
# 5. Selection cosmetic update
# ---------------------------------------------------------------------------
$ws.Range("C4").Select()
